$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("divisi") between jabatan (C) and awal_masuk (old D)
$ws.Range("D1").EntireColumn.Insert()

# Insert a new row 3 so we have two data rows
$ws.Range("A3").EntireRow.Insert()

# --- Header row ---
$ws.Range("A1").Value = "nama"
$ws.Range("B1").Value = "nomor_aplikasi"
$ws.Range("C1").Value = "jabatan"
$ws.Range("D1").Value = "divisi"
$ws.Range("E1").Value = "awal_masuk"
$ws.Range("F1").Value = "cuti"

# --- Row 2 ---
$ws.Range("A2").Value = "sss"
$ws.Range("B2").Value = 123456
$ws.Range("C2").Value = "'STAFF"
$ws.Range("D2").Value = "IT"
$ws.Range("E2").Value = "'2019-12-01"
$ws.Range("F2").Value = 0

# --- Row 3 ---
$ws.Range("A3").Value = "asd"
$ws.Range("B3").Value = 324324
$ws.Range("C3").Value = "'STAFF"
$ws.Range("D3").Value = "MARKETING"
$ws.Range("E3").Value = "'2019-12-02"
$ws.Range("F3").Value = 0

# --- Column widths ---
$ws.Columns("C").ColumnWidth = 20.3
$ws.Columns("D").ColumnWidth = 17.0

# --- Data validations ---
$ws.Range("C2:C3").Validation.Add(3, 1, 1, '"DIREKTUR UTAMA,DIREKTUR KEPATUHAN,DIREKTUR MARKETING,MANAGER,STAFF"')
$ws.Range("D2:D3").Validation.Add(3, 1, 1, '"BOD,MARKETING,IT"')

# --- Selection ---
$ws.Range("D5").Select()

Write-Host "done"
